$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows of username/password data (rows 3-7) ---
$ws.Range("A3").Value = "a@a"
$ws.Range("A4").Value = "b@b"
$ws.Range("A5").Value = "c@c"
$ws.Range("A6").Value = "z@z"
$ws.Range("A7").Value = "elisadiskg5@gmail.com"

$ws.Range("B3").Value = "a"
$ws.Range("B4").Value = "b"
$ws.Range("B5").Value = "c"
$ws.Range("B6").Value = "z"
$ws.Range("B7").Value = "star"

# --- Turn the email addresses in column A into mailto hyperlinks ---
# (this also creates the built-in "Hyperlink" cell style used by these cells)
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:a@a") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:b@b") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:c@c") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:z@z") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:elisadiskg5@gmail.com") | Out-Null

# --- Widen column A so the email addresses are fully visible ---
$ws.Columns("A").ColumnWidth = 22.375

# --- Page setup: portrait, paper size 9 (A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Final selection lands on B4 ---
$ws.Range("B4").Select() | Out-Null
